# Generate Report for Handback
# - Flip the "Ready for handoff" status to the handed-back status on both
#   locale sheets.
# - Stamp the handback datetime for each locale.
# - Add "Latest Target File" (F) / "Latest Handback File" (G) columns with
#   hyperlinked file names for row 2 and row 3 on both locale sheets.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$zhXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$srcMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ff8e9515857bd4c05f5e66a58bce17ead003f619/e2e/a.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb29afe583b30b9534ebc3e81d1eaefb6a5042c9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56801cbe40e999aff98551f9fd2620125ad6824c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

function Update-LocaleSheet {
    param($SheetName, $HandbackStamp, $XlfName, $XlfUrl)

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($row in 2, 3) {
        # Status -> handed back
        $ws.Cells.Item($row, 3).Value = $newStatus

        # Latest Handback DateTime
        $ws.Cells.Item($row, 8).Value = $HandbackStamp

        # F: Latest Target File
        $ws.Hyperlinks.Add($ws.Range("F$row"), $srcMdUrl, "", "", "a.md") | Out-Null

        # G: Latest Handback File
        $ws.Hyperlinks.Add($ws.Range("G$row"), $XlfUrl, "", "", $XlfName) | Out-Null
    }
}

Update-LocaleSheet "zh-cn" "2016-03-20 00:25:14" $zhXlf $zhXlfUrl
Update-LocaleSheet "de-de" "2016-03-20 00:25:19" $deXlf $deXlfUrl

# The "Ready for handoff" shared string is also surfaced on the Overview
# sheet (per-locale status columns B & C) - keep it in sync the same way
# the per-locale sheets were updated above.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
